$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "28.247.37"
Set-TextValue $ws.Range("E2") "  +3.07%  "
Set-TextValue $ws.Range("D3") "1.823.50"
Set-TextValue $ws.Range("E3") "  +1.58%  "
Set-TextValue $ws.Range("D4") "1.002"
Set-TextValue $ws.Range("E4") "  -0.23%  "
Set-TextValue $ws.Range("D5") "339.68"
Set-TextValue $ws.Range("E5") "  +0.65%  "
Set-TextValue $ws.Range("D6") "0.9966"
Set-TextValue $ws.Range("E6") "  -0.42%  "
Set-TextValue $ws.Range("D7") "0.3958"
Set-TextValue $ws.Range("E7") "  +4.23%  "
Set-TextValue $ws.Range("D8") "0.3501"
Set-TextValue $ws.Range("E8") "  +1.39%  "
Set-TextValue $ws.Range("E9") "  -0.61%  "
Set-TextValue $ws.Range("D10") "1.202"
Set-TextValue $ws.Range("E10") "  -0.02%  "
Set-TextValue $ws.Range("D11") "0.07596"
Set-TextValue $ws.Range("E11") "  +1.38%  "
Set-TextValue $ws.Range("D12") "0.9953"
Set-TextValue $ws.Range("E12") "  -0.64%  "
Set-TextValue $ws.Range("D13") "22.22"
Set-TextValue $ws.Range("E13") "  +0.74%  "
Set-TextValue $ws.Range("D14") "6.543"
Set-TextValue $ws.Range("E14") "  +1.13%  "
Set-TextValue $ws.Range("D15") "1.827.04"
Set-TextValue $ws.Range("E15") "  +1.67%  "
Set-TextValue $ws.Range("D16") "7.210"
Set-TextValue $ws.Range("E16") "  +2.15%  "
Set-TextValue $ws.Range("D17") "0.00001109"
Set-TextValue $ws.Range("E17") "  +0.81%  "
Set-TextValue $ws.Range("D18") "0.06709"
Set-TextValue $ws.Range("E18") "  +0.98%  "
Set-TextValue $ws.Range("D19") "85.51"
Set-TextValue $ws.Range("E19") "  +0.87%  "
Set-TextValue $ws.Range("D20") "0.9967"
Set-TextValue $ws.Range("E20") "  -0.37%  "
Set-TextValue $ws.Range("D21") "17.86"
Set-TextValue $ws.Range("E21") "  +3.10%  "
Set-TextValue $ws.Range("D22") "6.598"
Set-TextValue $ws.Range("E22") "  +1.23%  "
Set-TextValue $ws.Range("D23") "28.303.12"
Set-TextValue $ws.Range("E23") "  +3.34%  "
Set-TextValue $ws.Range("D24") "12.79"
Set-TextValue $ws.Range("E24") "  +2.14%  "
Set-TextValue $ws.Range("D25") "2.407"
Set-TextValue $ws.Range("E25") "  -1.07%  "
Set-TextValue $ws.Range("D26") "2.588"
Set-TextValue $ws.Range("E26") "  +1.14%  "
Set-TextValue $ws.Range("D27") "1.498"
Set-TextValue $ws.Range("E27") "  -0.35%  "
Set-TextValue $ws.Range("E28") "  +0.68%  "
Set-TextValue $ws.Range("D29") "154.66"
Set-TextValue $ws.Range("E29") "  +1.77%  "
Set-TextValue $ws.Range("D30") "2.033.04"
Set-TextValue $ws.Range("E30") "  +1.65%  "
Set-TextValue $ws.Range("D31") "135.94"
Set-TextValue $ws.Range("E31") "  +1.50%  "
Set-TextValue $ws.Range("D32") "6.211"
Set-TextValue $ws.Range("E32") "  +1.52%  "
Set-TextValue $ws.Range("D33") "4.034"
Set-TextValue $ws.Range("E33") "  -0.59%  "
Set-TextValue $ws.Range("D34") "0.08854"
Set-TextValue $ws.Range("E34") "  +1.95%  "
Set-TextValue $ws.Range("D35") "13.27"
Set-TextValue $ws.Range("E35") "  +0.18%  "
Set-TextValue $ws.Range("E36") "  +1.82%  "
Set-TextValue $ws.Range("D37") "0.7008"
Set-TextValue $ws.Range("E37") "  +1.57%  "
Set-TextValue $ws.Range("D38") "0.02443"
Set-TextValue $ws.Range("E38") "  +4.60%  "
Set-TextValue $ws.Range("D39") "0.06577"
Set-TextValue $ws.Range("E39") "  +3.24%  "
Set-TextValue $ws.Range("D40") "1.616"
Set-TextValue $ws.Range("E40") "  -2.70%  "
Set-TextValue $ws.Range("E42") "  -0.45%  "
Set-TextValue $ws.Range("D43") "8.511"
Set-TextValue $ws.Range("E43") "  -4.17%  "
Set-TextValue $ws.Range("D44") "14.65"
Set-TextValue $ws.Range("E44") "  +1.93%  "
Set-TextValue $ws.Range("D45") "0.6520"
Set-TextValue $ws.Range("E45") "  +1.28%  "
Set-TextValue $ws.Range("D46") "3.893"
Set-TextValue $ws.Range("E46") "  +0.64%  "
Set-TextValue $ws.Range("D47") "2.178"
Set-TextValue $ws.Range("E47") "  +2.42%  "
Set-TextValue $ws.Range("D48") "132.01"
Set-TextValue $ws.Range("E48") "  +1.47%  "
Set-TextValue $ws.Range("D49") "0.07223"
Set-TextValue $ws.Range("E49") "  +0.41%  "
Set-TextValue $ws.Range("D50") "80.54"
Set-TextValue $ws.Range("E50") "  +1.12%  "
Set-TextValue $ws.Range("D51") "1.254"
Set-TextValue $ws.Range("E51") "  +2.86%  "
